$d = $word.ActiveDocument
$vt = [char]11

# --- Change 1: "Test 3" paragraph --------------------------------------
# Merge "Test 3" + " " into a single run "Test 3 " (no visible text change).
$d.Content.Find.Execute("Test 3 ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "Test 3 ", 2) | Out-Null

# Merge the manual line break with the following "Regola con il nome vuoto "
# run (no visible text change).
$search1b = $vt + "Regola con il nome vuoto "
$d.Content.Find.Execute($search1b, $true, $false, $false, $false, $false, `
    $true, 1, $false, $search1b, 2) | Out-Null

# --- Change 2: "al nome della Regola" -----------------------------------
# Merge " al nome de" + "lla Regola nella finestra di dialogo che appare"
# into a single run (no visible text change, just a typo split fixed).
$search2 = " al nome della Regola nella finestra di dialogo che appare"
$d.Content.Find.Execute($search2, $true, $false, $false, $false, $false, `
    $true, 1, $false, $search2, 2) | Out-Null

# --- Change 3: replace the yellow-highlighted BUG note with the new ----
# "Errore nella definizione dei requisiti" sentence, and drop the
# highlighting.
$search3 = "BUG vi è la possibilità di avere una Regola con nome vuoto, questo bug è stato riscontrato nella fase di testing"
$replace3 = "Errore nella definizione dei requisiti: è possibile creare una Regola con il nome vuoto"
$d.Content.Find.Execute($search3, $true, $false, $false, $false, $false, `
    $true, 1, $false, $replace3, 2) | Out-Null

$r3 = $d.Content
$r3.Find.Execute($replace3, $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$r3.HighlightColorIndex = 0

# --- Reposition the "_GoBack" last-edit bookmark ------------------------
# Word re-drops its hidden "_GoBack" bookmark at the site of the most
# recent edit. The final edit in this revision lands in the "Come
# evidenziato..." paragraph, splitting it right after "...attraverso i".
try {
    $old = $d.Bookmarks("_GoBack")
    $old.Delete() | Out-Null
} catch {
}
$r4 = $d.Content
$r4.Find.Execute("attraverso i", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$bmRange = $d.Range($r4.End, $r4.End)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
